$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking scheme: correct-answer marks increased from 3 to 5
$ws.Range("B11").Value = 5

# Update total marks obtained for correct answers (22 right * 5 = 110)
$ws.Range("B12").Value = 110

# Update the correct/total marks summary text
$ws.Range("E12").Value = "110/140"
